# Generate Report for Handoff
#
# This script reflects that the localization-status report was regenerated:
# the row that used to describe "e94f9093-...md" (row 6) now describes
# "2257adb3-...md", and the row that used to describe "2257adb3-...md"
# (row 7) now describes "e94f9093-...md" -- together with refreshed
# status/timestamp/handoff-file data for both rows, across all three
# worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay($ws, [string]$addr, [string]$text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A6").Value = "2257adb3-2956-4548-b452-34068782e39a.md"
$wsOverview.Range("B6").Value = "e2e\2257adb3-2956-4548-b452-34068782e39a.md"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2017-02-17 06:51:17"

$wsOverview.Range("A7").Value = "e94f9093-06e5-4986-940f-c51c76a68213.md"
$wsOverview.Range("B7").Value = "e2e\e94f9093-06e5-4986-940f-c51c76a68213.md"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2017-02-17 06:58:30"

Set-HyperlinkDisplay $wsOverview '$B$6' "e2e\2257adb3-2956-4548-b452-34068782e39a.md"
Set-HyperlinkDisplay $wsOverview '$B$7' "e2e\e94f9093-06e5-4986-940f-c51c76a68213.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A6").Value = "2257adb3-2956-4548-b452-34068782e39a.md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("G6").Value = "2257adb3-2956-4548-b452-34068782e39a.5e93699a6f5169092210189c78559e9aecb8190e.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2017-02-17 06:50:59"

$wsZhCn.Range("A7").Value = "e94f9093-06e5-4986-940f-c51c76a68213.md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("G7").Value = "e94f9093-06e5-4986-940f-c51c76a68213.236da46c2eda8ff9775ef7d969c5502f9c7defe9.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2017-02-17 06:58:10"

Set-HyperlinkDisplay $wsZhCn '$A$6' "2257adb3-2956-4548-b452-34068782e39a.md"
Set-HyperlinkDisplay $wsZhCn '$A$7' "e94f9093-06e5-4986-940f-c51c76a68213.md"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A6").Value = "2257adb3-2956-4548-b452-34068782e39a.md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("G6").Value = "2257adb3-2956-4548-b452-34068782e39a.5e93699a6f5169092210189c78559e9aecb8190e.de-de.xlf"
$wsDeDe.Range("H6").Value = "2017-02-17 06:51:17"

$wsDeDe.Range("A7").Value = "e94f9093-06e5-4986-940f-c51c76a68213.md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("G7").Value = "e94f9093-06e5-4986-940f-c51c76a68213.236da46c2eda8ff9775ef7d969c5502f9c7defe9.de-de.xlf"
$wsDeDe.Range("H7").Value = "2017-02-17 06:58:30"

Set-HyperlinkDisplay $wsDeDe '$A$6' "2257adb3-2956-4548-b452-34068782e39a.md"
Set-HyperlinkDisplay $wsDeDe '$A$7' "e94f9093-06e5-4986-940f-c51c76a68213.md"
